$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1. Clear the old "Brokerage (%)" data values (header K1 stays put)
# -----------------------------------------------------------------
$ws.Range("K2:K4").ClearContents()

# -----------------------------------------------------------------
# 2. Make room for the new "Cedent/Reinsurer" + "Ceding Commission"
#    columns. Insert 3 blank columns right before column L; this
#    shifts the old L (Gross Premium), M (Share Premium) and N
#    (Ceding Comission) columns to O, P, Q, and everything after
#    them (Gross Brokerage Amount ... Processed) ends up at R..V,
#    which is exactly where the new layout needs them.
# -----------------------------------------------------------------
$ws.Columns.Item(12).Insert()
$ws.Columns.Item(12).Insert()
$ws.Columns.Item(12).Insert()

# -----------------------------------------------------------------
# 3. Re-write the header row for the six columns L..Q with the new
#    labels (this also overwrites the old Gross Premium / Share
#    Premium / Ceding Comission headers that were pushed into O/P/Q).
# -----------------------------------------------------------------
$ws.Range("L1").Value = "Cedent Gross Premium (100%)"
$ws.Range("M1").Value = "Cedent Share Premium"
$ws.Range("N1").Value = "Reinsurer Gross Premium (100%)"
$ws.Range("O1").Value = "Reinsurer Share Premium"
$ws.Range("P1").Value = "Ceding Commission (%)"
$ws.Range("Q1").Value = "Ceding Commission (Amount)"

# -----------------------------------------------------------------
# 4. Row 2 (existing debit note RI/25-26/GIFT/D08) - recompute the
#    premium split using the Cedent Rate (H) and Reinsurance Rate
#    (I) separately, and zero out the (now separated) ceding
#    commission / brokerage figures.
# -----------------------------------------------------------------
$ws.Range("L2").Value = 111169.17360000001
$ws.Range("M2").Value = 7781.8421520000011
$ws.Range("N2").Value = 111169.17360000001
$ws.Range("O2").Value = 7781.8421520000011
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 7781.8421520000011
$ws.Range("U2").Value = 7781.8421520000011
$ws.Range("V2").Value = "Yes"

# -----------------------------------------------------------------
# 5. Row 3 - new debit note RI/25-26/GIFT/D19 (Alibey Maldives /
#    Joali Maldives Resort placement), where the Reinsurance Rate
#    (7%) is lower than the Cedent Rate (10%).
# -----------------------------------------------------------------
$ws.Range("A3").Value = "RI/25-26/GIFT/D19"
$ws.Range("B3").Value = "04th November 2025"
$ws.Range("C3").Value = "Property All Risk "
$ws.Range("D3").Value = "Alibey Maldives Pvt Ltd and/or Joali Maldives Resort"
$ws.Range("E3").Value = "Solarelle Insurance Pvt. Ltd"
$ws.Range("F3").Value = "18-09-2025 To 17-09-2026"
$ws.Range("G3").Value = 183629631.02000001
$ws.Range("H3").Value = 0.1
$ws.Range("I3").Value = 0.07
$ws.Range("J3").Value = 13
$ws.Range("L3").Value = 183629.63102
$ws.Range("M3").Value = 23871.8520326
$ws.Range("N3").Value = 128540.74171400003
$ws.Range("O3").Value = 16710.296422820003
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 23871.8520326
$ws.Range("U3").Value = 16710.296422820003
$ws.Range("V3").Value = "Yes"

# -----------------------------------------------------------------
# 6. Row 4 - new debit note RI/25-26/GIFT/D20, same placement as
#    row 3, lower share.
# -----------------------------------------------------------------
$ws.Range("A4").Value = "RI/25-26/GIFT/D20"
$ws.Range("B4").Value = "04th November 2025"
$ws.Range("C4").Value = "Property All Risk "
$ws.Range("D4").Value = "Alibey Maldives Pvt Ltd and/or Joali Maldives Resort"
$ws.Range("E4").Value = "Solarelle Insurance Pvt. Ltd"
$ws.Range("F4").Value = "18-09-2025 To 17-09-2026"
$ws.Range("G4").Value = 10000000
$ws.Range("H4").Value = 0.2
$ws.Range("I4").Value = 0.18
$ws.Range("J4").Value = 15
$ws.Range("L4").Value = 20000
$ws.Range("M4").Value = 3000
$ws.Range("N4").Value = 18000
$ws.Range("O4").Value = 2700
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 3000
$ws.Range("U4").Value = 2700
$ws.Range("V4").Value = "Yes"

# -----------------------------------------------------------------
# 7. Resize the columns to fit their (new) content, then mirror the
#    final workbook's selection / scrolled view over to the new
#    columns.
# -----------------------------------------------------------------
$ws.Range("A1:Z4").Columns.AutoFit()
$ws.Range("V4").Select()
